$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D column to Text format to prevent Excel auto-converting numeric-looking strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.862.74'
$ws.Range("E2").Value = '  +1.02%  '

$ws.Range("D3").Value = '3.919.36'
$ws.Range("E3").Value = '  +2.22%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.47%  '

$ws.Range("D5").Value = '473.79'
$ws.Range("E5").Value = '  +5.70%  '

$ws.Range("D6").Value = '144.74'
$ws.Range("E6").Value = '  -2.53%  '

$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -1.81%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = '0.717'
$ws.Range("E9").Value = '  -3.31%  '

$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +5.19%  '

$ws.Range("D11").Value = '0.0000341'
$ws.Range("E11").Value = '  +5.13%  '

$ws.Range("D12").Value = '42.42'
$ws.Range("E12").Value = '  -3.21%  '

$ws.Range("D13").Value = '4.612.44'
$ws.Range("E13").Value = '  +3.90%  '

$ws.Range("D14").Value = '10.18'
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.989.74'
$ws.Range("E15").Value = '  +4.75%  '

$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").Value = '14.58'
$ws.Range("E16").Value = '  -3.00%  '

$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").Value = '19.66'
$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("D19").Value = '1.12'
$ws.Range("E19").Value = '  -2.87%  '

$ws.Range("D20").Value = '68.284.75'
$ws.Range("E20").Value = '  +1.61%  '

$ws.Range("D21").Value = '430.32'
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").Value = '3.26'
$ws.Range("E22").Value = '  +0.64%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '14.14'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("D24").Value = '86.55'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").Value = '3.60'
$ws.Range("E25").Value = '  +4.84%  '

$ws.Range("D26").Value = '37.86'
$ws.Range("E26").Value = '  +1.21%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '9.74'
$ws.Range("E27").Value = '  +1.54%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  +1.59%  '

$ws.Range("D29").Value = '713.83'
$ws.Range("E29").Value = '  -4.83%  '

$ws.Range("D30").Value = '13.09'
$ws.Range("E30").Value = '  -4.67%  '

$ws.Range("D31").Value = '0.125'
$ws.Range("E31").Value = '  -6.22%  '

$ws.Range("D32").Value = '2.79'
$ws.Range("E32").Value = '  +2.03%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '63.65'
$ws.Range("E33").Value = '  +10.02%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").Value = '41.30'
$ws.Range("E34").Value = '  -4.00%  '

$ws.Range("D35").Value = '0.0₃0821'
$ws.Range("E35").Value = '  +20.21%  '

$ws.Range("D36").Value = '0.148'
$ws.Range("E36").Value = '  -4.81%  '

$ws.Range("D37").Value = '0.997'
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("D38").Value = '5.25'
$ws.Range("E38").Value = '  -5.21%  '

$ws.Range("D39").Value = '0.0462'
$ws.Range("E39").Value = '  -2.72%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '30.70'
$ws.Range("E40").Value = '  +22.56%  '

$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").Value = '2.99'
$ws.Range("E41").Value = '  +3.99%  '

$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = '2.61'
$ws.Range("E42").Value = '  +5.62%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '2.88'
$ws.Range("E43").Value = '  +8.50%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.35%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.139'
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '0.327'
$ws.Range("E46").Value = '  -4.68%  '

$ws.Range("B47").Value = 'LidoDAOToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").Value = '  -2.96%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '2.12'
$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.17'
$ws.Range("E49").Value = '  -2.24%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '145.88'
$ws.Range("E50").Value = '  -0.45%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.87'
$ws.Range("E51").Value = '  -0.32%  '

# Reset D column style back to the default (Normal) so no stray formatting remains
$ws.Range("D2:D51").Style = "Normal"
